# Generate Report for Archive
#
# The three status-report rows whose "Source File Name" is
#   8b04af71-d0dc-45ff-846e-a6fb30c4252e
#   98279895-6a9d-4e67-a293-2da57e777709
#   cb513c5a-9a1e-4d44-8e7e-cfd332acf494
# are re-ordered (the 8b04af71 row moves from the top of the group to the
# bottom, the other two shift up) on every sheet: Overview, zh-cn, de-de.
# The 8b04af71 row keeps the "Ready for handoff" status while the two rows
# that shift into its old slot pick up "In Translation" instead.

$wb = $excel.ActiveWorkbook

# =======================================================================
# Sheet "Overview"  (columns: A=File Name, B=zh-cn status, C=de-de status,
#                    D=Latest Handoff Date)
# =======================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A8").Value = "98279895-6a9d-4e67-a293-2da57e777709.md"
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"

$wsOverview.Range("A9").Value = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md"
$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

$wsOverview.Range("A10").Value = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md"
$wsOverview.Range("B10").Value = "Ready for handoff"
$wsOverview.Range("C10").Value = "Ready for handoff"

# Re-sync the cached hyperlink display text for column A (rows 8-10 carry
# hyperlinks to the source-file commit on GitHub; the link target stays
# pinned to the row, only the cached display label follows the new value).
$i = 0
foreach ($link in $wsOverview.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 7) { $link.TextToDisplay = "98279895-6a9d-4e67-a293-2da57e777709.md" }
    if ($i -eq 8) { $link.TextToDisplay = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md" }
    if ($i -eq 9) { $link.TextToDisplay = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md" }
}

# =======================================================================
# Sheet "zh-cn"  (columns: A=Source File Name, B=File Extension,
#                 C=Status, D=Latest Handoff File, E=Latest Handoff Datetime)
# =======================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A8").Value = "98279895-6a9d-4e67-a293-2da57e777709.md"
$wsZh.Range("C8").Value = "In Translation"
$wsZh.Range("D8").Value = "98279895-6a9d-4e67-a293-2da57e777709.ceee6ed478cab9f4eb38bffce2162ec098cff7e3.zh-cn.xlf"

$wsZh.Range("A9").Value = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md"
$wsZh.Range("C9").Value = "In Translation"
$wsZh.Range("D9").Value = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.fa12590faf79ddb9993a07118a420541bb913fa6.zh-cn.xlf"

$wsZh.Range("A10").Value = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md"
$wsZh.Range("C10").Value = "Ready for handoff"
$wsZh.Range("D10").Value = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.1256191c502f5bdcc482b405b385b12dd89fdd69.zh-cn.xlf"

$i = 0
foreach ($link in $wsZh.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 27) { $link.TextToDisplay = "98279895-6a9d-4e67-a293-2da57e777709.md" }
    if ($i -eq 29) { $link.TextToDisplay = "98279895-6a9d-4e67-a293-2da57e777709.ceee6ed478cab9f4eb38bffce2162ec098cff7e3.zh-cn.xlf" }
    if ($i -eq 30) { $link.TextToDisplay = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md" }
    if ($i -eq 32) { $link.TextToDisplay = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.fa12590faf79ddb9993a07118a420541bb913fa6.zh-cn.xlf" }
    if ($i -eq 33) { $link.TextToDisplay = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md" }
    if ($i -eq 35) { $link.TextToDisplay = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.1256191c502f5bdcc482b405b385b12dd89fdd69.zh-cn.xlf" }
}

# =======================================================================
# Sheet "de-de"  (same layout as zh-cn)
# =======================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A8").Value = "98279895-6a9d-4e67-a293-2da57e777709.md"
$wsDe.Range("C8").Value = "In Translation"
$wsDe.Range("D8").Value = "98279895-6a9d-4e67-a293-2da57e777709.ceee6ed478cab9f4eb38bffce2162ec098cff7e3.de-de.xlf"

$wsDe.Range("A9").Value = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md"
$wsDe.Range("C9").Value = "In Translation"
$wsDe.Range("D9").Value = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.fa12590faf79ddb9993a07118a420541bb913fa6.de-de.xlf"

$wsDe.Range("A10").Value = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md"
$wsDe.Range("C10").Value = "Ready for handoff"
$wsDe.Range("D10").Value = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.1256191c502f5bdcc482b405b385b12dd89fdd69.de-de.xlf"

$i = 0
foreach ($link in $wsDe.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 27) { $link.TextToDisplay = "98279895-6a9d-4e67-a293-2da57e777709.md" }
    if ($i -eq 29) { $link.TextToDisplay = "98279895-6a9d-4e67-a293-2da57e777709.ceee6ed478cab9f4eb38bffce2162ec098cff7e3.de-de.xlf" }
    if ($i -eq 30) { $link.TextToDisplay = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.md" }
    if ($i -eq 32) { $link.TextToDisplay = "cb513c5a-9a1e-4d44-8e7e-cfd332acf494.fa12590faf79ddb9993a07118a420541bb913fa6.de-de.xlf" }
    if ($i -eq 33) { $link.TextToDisplay = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.md" }
    if ($i -eq 35) { $link.TextToDisplay = "8b04af71-d0dc-45ff-846e-a6fb30c4252e.1256191c502f5bdcc482b405b385b12dd89fdd69.de-de.xlf" }
}
